# Auto-generated edit script: updates cryptos list (price/volume columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "65.653.60"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.675.88"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.39"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.69"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +6.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.129"
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.37"
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "3.155.28"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "65.499.26"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "2.667.30"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.81"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.69"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.62"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000111"
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -5.46%  "
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "531.53"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.47"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.424"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.51"
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "157.94"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "164.35"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.16"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0610"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.643"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("E49").Value = "  +15.21%  "
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.11"
$ws.Range("E51").Value = "  -5.53%  "
